$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 was a duplicate "Ore and mineral industries" row accidentally left in
# the sheet (row 13 already holds that sector). Delete it so everything below
# shifts up by one row, fixing the table / sheet dimensions accordingly.
$ws.Rows.Item(14).Delete()

# Update the active selection to match the saved state in the fixed file.
$ws.Range("A14:XFD14").Select()
